# Commit: "unify the conception of DataNode, DataTable, Entity."
# The only substantive/object-model-visible change in this revision is the
# rename of the sole worksheet from "Property1" to "DataNode" (the rest of
# the diff - fileVersion/rupBuild bumps, the absPath machine path, and the
# various xr/xr2/xr16/x16r2 revision-tracking GUIDs - are artifacts Excel
# stamps on save and are not reachable/meaningful through the object model).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DataNode"
